$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of values to write across A1:BI1 (61 columns)
$values = @(0,1,2,4,5,6,8,9,10,11,12,14,15,17,19,20,21,22,23,25,26,27,29,31,32,34,36,37,38,39,40,41,42,43,44,45,46,47,48,49,51,52,53,54,55,57,58,59,60,62,63,64,65,66,67,68,69,70,71,73,75)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $values[$i]
}
